# Add data for a second UQAM site (UQAM_2) so the app can load data
# from multiple sites, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: first new record for UQAM_2
$ws.Range("A3").Value = "UQAM_2"
$ws.Range("B3").Value = -75
$ws.Range("C3").Value = 46.531700000000001
$ws.Range("D3").Value = -72.6447

# Row 4: second new record for UQAM_2
$ws.Range("A4").Value = "UQAM_2"
$ws.Range("B4").Value = -75
$ws.Range("C4").Value = 46.527971999999998
$ws.Range("D4").Value = -72.646249999999995

# Leave the selection where Excel would naturally end up after entering
# this data (matches the saved cursor position in the workbook).
$ws.Range("D7").Select()
